# fix(gui) step 1 and 2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: bump the date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the computed price column (D33:D37)
$ws.Range("D33").Value = 9431.802
$ws.Range("D34").Value = 11160.692
$ws.Range("D35").Value = 13785.194
$ws.Range("D36").Value = 15376.297
$ws.Range("D37").Value = 15901.197
